$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '29.909.76'
$ws.Range("E2").Value = '  +0.45%  '

# Row 3
$ws.Range("D3").Value = '1.895.37'
$ws.Range("E3").Value = '  +0.27%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.10%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7806'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.28%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '243.95'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.69%  '

# Row 7
$ws.Range("E7").Value = '  +0.05%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3131'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -1.16%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '25.70'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +1.37%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07336'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +4.51%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08090'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.56%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7719'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.90%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.500'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +4.27%  '

# Row 14
$ws.Range("D14").Value = '1.899.97'
$ws.Range("E14").Value = '  +0.55%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '93.90'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +2.06%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.224'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +5.43%  '

# Row 17
$ws.Range("D17").Value = '29.838.25'
$ws.Range("E17").Value = '  +0.20%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.96'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.97%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '247.44'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +1.88%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007823'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +1.69%  '

# Row 21
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9996'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.03%  '

# Row 22
$ws.Range("B22").Value = 'Chainlink'
$ws.Range("C22").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.100'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -1.02%  '

# Row 23
$ws.Range("B23").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C23").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D23").Value = '2.116.82'
$ws.Range("E23").Value = '  -0.89%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.001'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.10%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1590'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -3.54%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.448'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +1.69%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '163.76'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -1.16%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.73'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.38%  '

# Row 29
$ws.Range("E29").Value = '  -1.07%  '

# Row 30
$ws.Range("E30").Value = '  +3.04%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.545'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.75%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.488'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +1.42%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05572'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.74%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.061'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.65%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.239'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -1.68%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7529'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +2.35%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.005'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.10%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.684'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +1.57%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01934'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +1.43%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.798'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +1.05%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4471'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +1.77%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '74.26'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +2.72%  '

# Row 43
$ws.Range("D43").Value = '1.112.26'
$ws.Range("E43").Value = '  +8.99%  '

# Row 44
$ws.Range("E44").Value = '  +2.83%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8512'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +1.57%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.0000'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.07%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.887'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +1.34%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '102.52'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.34%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.521'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +1.60%  '

# Row 50
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.777'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.05%  '

# Row 51
$ws.Range("B51").Value = 'SynthetixNetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.045'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +4.84%  '
